$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.279.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.591.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.81%  "

$ws.Range("E6").Value = "  +0.86%  "

$ws.Range("E8").Value = "  +0.28%  "

$ws.Range("E9").Value = "  -0.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0850"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.813.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.618.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.77%  "

$ws.Range("E14").Value = "  +0.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.522"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.49%  "

$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.280.42"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "213.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.94%  "

$ws.Range("E21").Value = "  -0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.95%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.12%  "

$ws.Range("E24").Value = "  -1.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.32%  "

$ws.Range("E26").Value = "  -0.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.88%  "

$ws.Range("E28").Value = "  -0.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0500"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.85%  "

$ws.Range("E31").Value = "  +1.17%  "

$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("E33").Value = "  +0.50%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.340.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.02%  "

$ws.Range("E35").Value = "  -0.91%  "

$ws.Range("E36").Value = "  -0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.594"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.06%  "

$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.818"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.28%  "

$ws.Range("E40").Value = "  +4.36%  "

$ws.Range("E41").Value = "  -0.23%  "

$ws.Range("E42").Value = "  -0.23%  "

$ws.Range("E43").Value = "  +0.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.764"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.30%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.56%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.724.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.03%  "

$ws.Range("E48").Value = "  -1.19%  "

$ws.Range("E49").Value = "  -3.58%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0980"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.26%  "

$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0504"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.40%  "
